$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row heights: rows 2-8 go from 15.5 to 15.75 (matching the default/other rows)
foreach ($r in 2..8) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# Selection: was a single cell (F6); now the whole data range A1:D20 is selected.
$ws.Range("A1:D20").Select()
